$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.823.55'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.753.71'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.06'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5094'
$ws.Range('E7').Value = '  +3.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.02'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2704'
$ws.Range('E9').Value = '  +10.35%  '
$ws.Range('E10').Value = '  +4.22%  '
$ws.Range('D11').Value = '1.741.73'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06919'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.60'
$ws.Range('E13').Value = '  +5.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6278'
$ws.Range('E14').Value = '  +8.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '78.70'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.499'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9995'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '25.838.68'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.72'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006739'
$ws.Range('E21').Value = '  +4.35%  '
$ws.Range('D22').Value = '1.973.67'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  +2.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.263'
$ws.Range('E24').Value = '  +4.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.184'
$ws.Range('E25').Value = '  +3.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.79'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.34'
$ws.Range('E27').Value = '  +5.42%  '
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.799'
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.86'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08274'
$ws.Range('E31').Value = '  +2.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.739'
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.432'
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04415'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.642'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.006'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6058'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.699'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.961'
$ws.Range('E39').Value = '  -4.38%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01565'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.36'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3867'
$ws.Range('E43').Value = '  +2.85%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7507'
$ws.Range('E44').Value = '  -3.47%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.877'
$ws.Range('E45').Value = '  -5.89%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05512'
$ws.Range('E46').Value = '  +7.67%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1097'
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.976'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.29'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.97'
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.40%  '

Write-Host "Updated cryptos list"
